# update code import data / get data from user
# Re-imports the latest material list prices from the user-supplied data,
# switches the sheet's font from Calibri to Arial, and leaves the selection
# on the last edited cell (matching the user's on-screen workflow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated list prices (mat_price column) ---------------------------------
$ws.Range("B2").Value  = 154000   # nabshi_metr
$ws.Range("B3").Value  = 302000   # f_47_metr
$ws.Range("B4").Value  = 219000   # u_36_metr
$ws.Range("B5").Value  = 960000   # panel_metr
$ws.Range("B6").Value  = 3050     # pich_panel_dane_adad
$ws.Range("B7").Value  = 2900     # pich_saze_dane_adad
$ws.Range("B8").Value  = 50000    # mikh_chashni_dane_adad
$ws.Range("B9").Value  = 17800    # ht_90_adad
$ws.Range("B10").Value = 37000    # beraket_adad
$ws.Range("B11").Value = 33500    # clips_adad
$ws.Range("B12").Value = 29400    # rabet_w_adad
# B13 (rolplak_adad) is unchanged at 40000

# --- Switch the workbook font from Calibri to Arial --------------------------
$ws.Range("A1:B13").Font.Name = "Arial"

# --- Restore the on-screen view: scrolled down with the last cell selected --
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B13").Select()
